# Swap the Species-related data (ASV_ID, Species_name, Common_name, Category)
# between rows 42 and 43, leaving ASV_sum (E) and ASV_rank (F) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture row 42 values (columns A-D)
$a42 = $ws.Range("A42").Value()
$b42 = $ws.Range("B42").Value()
$c42 = $ws.Range("C42").Value()
$d42 = $ws.Range("D42").Value()

# Capture row 43 values (columns A-D)
$a43 = $ws.Range("A43").Value()
$b43 = $ws.Range("B43").Value()
$c43 = $ws.Range("C43").Value()
$d43 = $ws.Range("D43").Value()

# Write row 43's original values into row 42
$ws.Range("A42").Value = $a43
$ws.Range("B42").Value = $b43
$ws.Range("C42").Value = $c43
$ws.Range("D42").Value = $d43

# Write row 42's original values into row 43
$ws.Range("A43").Value = $a42
$ws.Range("B43").Value = $b42
$ws.Range("C43").Value = $c42
$ws.Range("D43").Value = $d42
